$p = $ppt.ActivePresentation
$p.Slides.Item(12).Delete()
